$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers (shorter labels so the exported sheet is narrower)
$ws.Range("F1").Value = "Memory tot.size"
$ws.Range("H1").Value = "Memory P/Ns"
$ws.Range("M1").Value = "HDD slot pop."
$ws.Range("N1").Value = "PSU P/Ns"

# Shrink the columns to match the new, shorter header text
$ws.Columns.Item(6).ColumnWidth = 14.83
$ws.Columns.Item(8).ColumnWidth = 10.83
$ws.Columns.Item(13).ColumnWidth = 12.83
$ws.Columns.Item(14).ColumnWidth = 7.83
